# Add a header row to the "lessonList" sheet: Subject / Teacher / StudentGroup
# above the existing lesson rows (A=#, B=Subject, C=Teacher, D=StudentGroup,
# E/F=flags on the first data row only).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lessonList")

# Push all existing data down by one row, freeing up row 1 for headers.
$ws.Rows.Item(1).Insert()

# The data rows use style index 2 (wrapped Courier New) for columns B:D.
# Copy that formatting onto the new header row so it matches the sheet's
# existing look instead of falling back to the default style.
$ws.Range("B2:D2").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the header labels (column A/E/F stay blank on the header row).
$ws.Range("B1").Value = "Subject"
$ws.Range("C1").Value = "Teacher"
$ws.Range("D1").Value = "StudentGroup"

# Match the taller header row and the widened StudentGroup column.
$ws.Rows.Item(1).RowHeight = 23.85
$ws.Columns.Item(4).ColumnWidth = 27.6466666666667

# Leave the selection on the new header cell.
$ws.Range("B1").Select()
